$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Row 1 (header): replace numeric sequence 1..8 with P1..P9 labels, add column J ---
$ws.Range("B1").Value = "P1"
$ws.Range("C1").Value = "P2"
$ws.Range("D1").Value = "P3"
$ws.Range("E1").Value = "P4"
$ws.Range("F1").Value = "P5"
$ws.Range("G1").Value = "P6"
$ws.Range("H1").Value = "P7"
$ws.Range("I1").Value = "P8"
$ws.Range("J1").Value = "P9"

# --- Row 2 (data): G2 changes answer to "NA"; shift the "na grade do ppge" note into
#     the new column J, backfilling I2 with "não" ---
$ws.Range("G2").Value = "NA"
$ws.Range("I2").Value = "não"
$ws.Range("J2").Value = "na grade do ppge"

# --- Header styling: black font, centered horizontally & vertically. Stage the combined
#     format on a scratch cell and paste it in one shot so the header row ends up on a
#     single new style record instead of one per property. ---
$stage = $ws.Range("Z1")
$stage.Font.Color = 0
$stage.HorizontalAlignment = -4108
$stage.VerticalAlignment = -4108
$stage.Copy()
$ws.Range("A1:J1").PasteSpecial(-4122)
$stage.Clear()
$excel.CutCopyMode = 0

# --- Column widths (approximate the auto-fit widths Excel computed for the new data) ---
$ws.Columns.Item(1).ColumnWidth = 6.498697916666667
$ws.Columns.Item(2).ColumnWidth = 7.498697916666667
$ws.Columns.Item(3).ColumnWidth = 19.166666666666668
$ws.Columns.Item(4).ColumnWidth = 2.9986979166666665
$ws.Columns.Item(5).ColumnWidth = 2.6666666666666665
$ws.Columns.Item(6).ColumnWidth = 2.9986979166666665
$ws.Columns.Item(7).ColumnWidth = 2.6666666666666665
$ws.Columns.Item(8).ColumnWidth = 3.1666666666666665
$ws.Columns.Item(9).ColumnWidth = 3.1666666666666665
$ws.Columns.Item(10).ColumnWidth = 13.498697916666666

# --- Page setup tweak (A4 portrait) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moves to B2 ---
$ws.Range("B2").Select() | Out-Null
